# Auto-generated Excel COM-interop script
# Applies scheduled market-price/profit data refresh to the Kujata_Profits workbook
# across all Leve-crafting job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2912.3718
$ws.Range("I15").Value = 2912.3718
$ws.Range("K15").Value = 8737.115399999999
$ws.Range("M15").Value = -8568.115399999999
# Row 17
$ws.Range("H17").Value = 2014.4445
$ws.Range("J17").Value = 2014.4445
$ws.Range("L17").Value = 6043.333500000001
$ws.Range("N17").Value = -6379.333500000001
# Row 38
$ws.Range("H38").Value = 2072.75
$ws.Range("I38").Value = 20.5
$ws.Range("J38").Value = 4125
$ws.Range("K38").Value = 61.5
$ws.Range("L38").Value = 12375
$ws.Range("M38").Value = 310.5
$ws.Range("N38").Value = -13119
# Row 92
$ws.Range("H92").Value = 1983.8
$ws.Range("I92").Value = 2128.923
$ws.Range("K92").Value = 2128.923
$ws.Range("M92").Value = -880.9229999999998
# Row 106
$ws.Range("H106").Value = 11007.462
$ws.Range("I106").Value = 12325.909
$ws.Range("J106").Value = 3756
$ws.Range("K106").Value = 12325.909
$ws.Range("L106").Value = 3756
$ws.Range("M106").Value = -11694.909
$ws.Range("N106").Value = -5018
# Row 111
$ws.Range("H111").Value = 2569.9333
$ws.Range("I111").Value = 2324.3333
$ws.Range("J111").Value = 2938.3333
$ws.Range("K111").Value = 6972.999899999999
$ws.Range("L111").Value = 8814.999899999999
$ws.Range("M111").Value = -3905.999899999999
$ws.Range("N111").Value = -14948.9999
# Row 116
$ws.Range("H116").Value = 2229.6365
$ws.Range("I116").Value = 1669.3334
$ws.Range("K116").Value = 1669.3334
$ws.Range("M116").Value = 1772.6666
# Row 129
$ws.Range("H129").Value = 461.83334
$ws.Range("I129").Value = 336.05264
$ws.Range("J129").Value = 939.8
$ws.Range("K129").Value = 1008.15792
$ws.Range("L129").Value = 2819.4
$ws.Range("M129").Value = 3991.84208
$ws.Range("N129").Value = -12819.4
# Row 135
$ws.Range("H135").Value = 27027582
$ws.Range("I135").Value = 250.125
$ws.Range("J135").Value = 200002510
$ws.Range("K135").Value = 2251.125
$ws.Range("L135").Value = 1800022590
$ws.Range("M135").Value = 283.875
$ws.Range("N135").Value = -1800027660
# Row 137
$ws.Range("H137").Value = 1182.8823
$ws.Range("I137").Value = 885.56757
$ws.Range("J137").Value = 1537.742
$ws.Range("K137").Value = 2656.70271
$ws.Range("L137").Value = 4613.226
$ws.Range("M137").Value = -106.70271
$ws.Range("N137").Value = -9713.225999999999
# Row 138
$ws.Range("H138").Value = 1333.55
$ws.Range("I138").Value = 709.36365
$ws.Range("J138").Value = 1640.9851
$ws.Range("K138").Value = 2128.09095
$ws.Range("L138").Value = 4922.955300000001
$ws.Range("M138").Value = 3011.90905
$ws.Range("N138").Value = -15202.9553
# Row 141
$ws.Range("H141").Value = 735.625
$ws.Range("I141").Value = 651.3333
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 1953.9999
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 3226.0001
$ws.Range("N141").Value = -16360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4804.2046
$ws.Range("I32").Value = 4406.369
$ws.Range("J32").Value = 6240.8335
$ws.Range("K32").Value = 4406.369
$ws.Range("L32").Value = 6240.8335
$ws.Range("M32").Value = -4119.369
$ws.Range("N32").Value = -6814.8335
# Row 61
$ws.Range("H61").Value = 43479132
$ws.Range("I61").Value = 45455320
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 45455320
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -45455108
$ws.Range("N61").Value = -3424
# Row 74
$ws.Range("H74").Value = 1254.6666
$ws.Range("I74").Value = 977.9231
$ws.Range("K74").Value = 977.9231
$ws.Range("M74").Value = -103.9231
# Row 77
$ws.Range("H77").Value = 1254.6666
$ws.Range("I77").Value = 977.9231
$ws.Range("K77").Value = 4889.6155
$ws.Range("M77").Value = -521.6154999999999
# Row 132
$ws.Range("H132").Value = 1819.4103
$ws.Range("I132").Value = 1568.2
$ws.Range("J132").Value = 2268
$ws.Range("K132").Value = 4704.6
$ws.Range("L132").Value = 6804
$ws.Range("M132").Value = -2174.6
$ws.Range("N132").Value = -11864
# Row 136
$ws.Range("H136").Value = 43479132
$ws.Range("I136").Value = 45455320
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 136365960
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -136363410
$ws.Range("N136").Value = -14100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2660
$ws.Range("I86").Value = 3179.4666
$ws.Range("K86").Value = 3179.4666
$ws.Range("M86").Value = -2056.4666
# Row 89
$ws.Range("H89").Value = 2660
$ws.Range("I89").Value = 3179.4666
$ws.Range("K89").Value = 15897.333
$ws.Range("M89").Value = -10281.333
# Row 107
$ws.Range("H107").Value = 806.7727
$ws.Range("I107").Value = 635.08105
$ws.Range("K107").Value = 635.08105
$ws.Range("M107").Value = 1284.91895

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2126.7917
$ws.Range("I31").Value = 2243.15
$ws.Range("J31").Value = 1545
$ws.Range("K31").Value = 2243.15
$ws.Range("L31").Value = 1545
$ws.Range("M31").Value = -1948.15
$ws.Range("N31").Value = -2135
# Row 34
$ws.Range("H34").Value = 2126.7917
$ws.Range("I34").Value = 2243.15
$ws.Range("J34").Value = 1545
$ws.Range("K34").Value = 2243.15
$ws.Range("L34").Value = 1545
$ws.Range("M34").Value = -2041.15
$ws.Range("N34").Value = -1949
# Row 58
$ws.Range("H58").Value = 937.89795
$ws.Range("I58").Value = 876.2683
$ws.Range("J58").Value = 1253.75
$ws.Range("K58").Value = 876.2683
$ws.Range("L58").Value = 1253.75
$ws.Range("M58").Value = -673.2683
$ws.Range("N58").Value = -1659.75
# Row 94
$ws.Range("H94").Value = 913.15
$ws.Range("I94").Value = 957.7143
$ws.Range("K94").Value = 957.7143
$ws.Range("M94").Value = -506.7143
# Row 107
$ws.Range("H107").Value = 548.5
$ws.Range("I107").Value = 428
$ws.Range("J107").Value = 669
$ws.Range("K107").Value = 428
$ws.Range("L107").Value = 669
$ws.Range("M107").Value = 1492
$ws.Range("N107").Value = -4509
# Row 122
$ws.Range("H122").Value = 1412.5
$ws.Range("I122").Value = 1412.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4237.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1787.5
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 2738.2173
$ws.Range("I132").Value = 2512.6428
$ws.Range("J132").Value = 3089.111
$ws.Range("K132").Value = 7537.928400000001
$ws.Range("L132").Value = 9267.332999999999
$ws.Range("M132").Value = -5007.928400000001
$ws.Range("N132").Value = -14327.333
# Row 134
$ws.Range("H134").Value = 20834546
$ws.Range("I134").Value = 1129.7059
$ws.Range("K134").Value = 3389.1177
$ws.Range("M134").Value = -854.1176999999998
# Row 136
$ws.Range("H136").Value = 937.89795
$ws.Range("I136").Value = 876.2683
$ws.Range("J136").Value = 1253.75
$ws.Range("K136").Value = 2628.8049
$ws.Range("L136").Value = 3761.25
$ws.Range("M136").Value = -78.80490000000009
$ws.Range("N136").Value = -8861.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 940.24
$ws.Range("J122").Value = 1065
$ws.Range("L122").Value = 9585
$ws.Range("N122").Value = -14485
# Row 131
$ws.Range("H131").Value = 25003448
$ws.Range("J131").Value = 5053.769
$ws.Range("L131").Value = 15161.307
$ws.Range("N131").Value = -25241.307

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 86
$ws.Range("H86").Value = 30165.5
$ws.Range("J86").Value = 30165.5
$ws.Range("L86").Value = 30165.5
$ws.Range("N86").Value = -32537.5
# Row 89
$ws.Range("H89").Value = 30165.5
$ws.Range("J89").Value = 30165.5
$ws.Range("L89").Value = 90496.5
$ws.Range("N89").Value = -102352.5
# Row 132
$ws.Range("H132").Value = 2876.5
$ws.Range("I132").Value = 3112.3076
$ws.Range("K132").Value = 9336.9228
$ws.Range("M132").Value = -6806.9228

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1577.8096
$ws.Range("I7").Value = 1536.7059
$ws.Range("K7").Value = 1536.7059
$ws.Range("M7").Value = -1424.7059
# Row 22
$ws.Range("H22").Value = 654.8
$ws.Range("I22").Value = 549.2
$ws.Range("K22").Value = 549.2
$ws.Range("M22").Value = -254.2
# Row 27
$ws.Range("H27").Value = 654.8
$ws.Range("I27").Value = 549.2
$ws.Range("K27").Value = 549.2
$ws.Range("M27").Value = -442.2
# Row 40
$ws.Range("H40").Value = 2503.16
$ws.Range("I40").Value = 1844.7142
$ws.Range("K40").Value = 1844.7142
$ws.Range("M40").Value = -1708.7142
# Row 122
$ws.Range("H122").Value = 11809174
$ws.Range("I122").Value = 20239896
$ws.Range("K122").Value = 60719688
$ws.Range("M122").Value = -60717238
# Row 126
$ws.Range("H126").Value = 1577.8096
$ws.Range("I126").Value = 1536.7059
$ws.Range("K126").Value = 4610.1177
$ws.Range("M126").Value = -2140.1177
# Row 132
$ws.Range("H132").Value = 18846.914
$ws.Range("I132").Value = 1100.5128
$ws.Range("J132").Value = 55273.74
$ws.Range("K132").Value = 3301.5384
$ws.Range("L132").Value = 165821.22
$ws.Range("M132").Value = -771.5383999999999
$ws.Range("N132").Value = -170881.22

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3326.4055
$ws.Range("I132").Value = 3529.5217
$ws.Range("K132").Value = 10588.5651
$ws.Range("M132").Value = -8058.5651
# Row 136
$ws.Range("H136").Value = 476.6
$ws.Range("I136").Value = 412.76315
$ws.Range("K136").Value = 1238.28945
$ws.Range("M136").Value = 1311.71055

